# Update "Inscritos" (column E) counts on the "Inscricoes" worksheet
# to reflect the latest enrollment numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E7").Value = 30
$ws.Range("E11").Value = 17
$ws.Range("E16").Value = 309
$ws.Range("E18").Value = 95
